$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("T2").Value = 0.1426048780487805
$ws.Range("V2").Value = 0.0002220611491829204
$ws.Range("Z2").Value = -0.2234193421634373
$ws.Range("AB2").Value = -1006.116301683182
$ws.Range("AC2").Value = "umolO2/min/m2"
$ws.Range("AD2").Value = -1006.116301683182

# Row 3
$ws.Range("T3").Value = 0.1423414634146342
$ws.Range("V3").Value = 0.0002565454225970831
$ws.Range("Z3").Value = -0.1714821897008507
$ws.Range("AB3").Value = -668.4281791695488
$ws.Range("AC3").Value = "umolO2/min/m2"
$ws.Range("AD3").Value = -668.4281791695488

# Row 4
$ws.Range("T4").Value = 0.1477268292682927
$ws.Range("V4").Value = 0.0002491214197856265
$ws.Range("Z4").Value = -0.2280694701114503
$ws.Range("AB4").Value = -915.4952244078742
$ws.Range("AC4").Value = "umolO2/min/m2"
$ws.Range("AD4").Value = -915.4952244078742

# Row 5
$ws.Range("T5").Value = 0.1470341463414634
$ws.Range("V5").Value = 0.0002082674398172554
$ws.Range("Z5").Value = -0.2805450244660898
$ws.Range("AB5").Value = -1347.042171893286
$ws.Range("AC5").Value = "umolO2/min/m2"
$ws.Range("AD5").Value = -1347.042171893286

# Row 6
$ws.Range("T6").Value = 0.1477560975609756
$ws.Range("V6").Value = 0.0001607142857142857
$ws.Range("Z6").Value = -0.2022093500179011
$ws.Range("AB6").Value = -1258.191511222495
$ws.Range("AC6").Value = "umolO2/min/m2"
$ws.Range("AD6").Value = -1258.191511222495

# Row 7
$ws.Range("T7").Value = 0.1455317073170732
$ws.Range("V7").Value = 0.0002247627833421192
$ws.Range("Z7").Value = -0.1304589500380652
$ws.Range("AB7").Value = -580.4295003745754
$ws.Range("AC7").Value = "umolO2/min/m2"
$ws.Range("AD7").Value = -580.4295003745754

# Row 8
$ws.Range("T8").Value = 0.1544
$ws.Range("V8").Value = 0
$ws.Range("Z8").Value = 0.0002114954780799225
$ws.Range("AB8").Value = "Inf"
$ws.Range("AC8").Value = "umolO2/min/m2"
$ws.Range("AD8").Value = "Inf"

# Row 9
$ws.Range("T9").Value = 0.1426048780487805
$ws.Range("V9").Value = 0.0002220611491829204
$ws.Range("Z9").Value = 0.1914724901640171
$ws.Range("AB9").Value = 862.2511901273361
$ws.Range("AC9").Value = "umolO2/min/m2"
$ws.Range("AD9").Value = 862.2511901273361

# Row 10
$ws.Range("T10").Value = 0.1423414634146342
$ws.Range("V10").Value = 0.0002565454225970831
$ws.Range("Z10").Value = 18.07288897658993
$ws.Range("AB10").Value = 70447.13093546113
$ws.Range("AC10").Value = "umolO2/min/m2"
$ws.Range("AD10").Value = 70447.13093546113

# Row 11
$ws.Range("T11").Value = 0.1477268292682927
$ws.Range("V11").Value = 0.0002491214197856265
$ws.Range("Z11").Value = 0.2335392322475941
$ws.Range("AB11").Value = 937.4514341181856
$ws.Range("AC11").Value = "umolO2/min/m2"
$ws.Range("AD11").Value = 937.4514341181856

# Row 12
$ws.Range("T12").Value = 0.1470341463414634
$ws.Range("V12").Value = 0.0002082674398172554
$ws.Range("Z12").Value = 0.2918591234658028
$ws.Range("AB12").Value = 1401.367029440104
$ws.Range("AC12").Value = "umolO2/min/m2"
$ws.Range("AD12").Value = 1401.367029440104

# Row 13
$ws.Range("T13").Value = 0.1477560975609756
$ws.Range("V13").Value = 0.0001607142857142857
$ws.Range("Z13").Value = 0.2296609123500524
$ws.Range("AB13").Value = 1429.001232400326
$ws.Range("AC13").Value = "umolO2/min/m2"
$ws.Range("AD13").Value = 1429.001232400326

# Row 14
$ws.Range("T14").Value = 0.1455317073170732
$ws.Range("V14").Value = 0.0002247627833421192
$ws.Range("Z14").Value = 0.07747301022218814
$ws.Range("AB14").Value = 344.687893031934
$ws.Range("AC14").Value = "umolO2/min/m2"
$ws.Range("AD14").Value = 344.687893031934

# Row 15
$ws.Range("T15").Value = 0.1544
$ws.Range("V15").Value = 0
$ws.Range("Z15").Value = -0.0005546210473080581
$ws.Range("AB15").Value = "-Inf"
$ws.Range("AC15").Value = "umolO2/min/m2"
$ws.Range("AD15").Value = "-Inf"

